# Apply "Major Refactor to add categories to operators" edit to sampleMath.xlsx
# - Adds CEIL/CEILING, PI, MOD, FACT function demo rows
# - Re-shuffles the existing COMBIN/PERMUT/GCD/LCM rows down to make room
# - Updates the related named ranges
# - Cosmetic bits: column D width, selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-lay-out the C/D "function demo" column.
#    Original:  C10/D10=COMBIN  C11/D11=PERMUT  C13/D13=GCD  C14/D14=LCM
#    Target  :  C9/D9=CEIL(new) C11/D11=COMBIN  C12/D12=PERMUT
#               C14/D14=GCD     C15/D15=LCM     C17/D17=FACT(new)
#               C18/D18=PI(new) C19/D19=MOD(new)
# Work from the bottom up so we never overwrite a value before it has been
# relocated.
# ---------------------------------------------------------------------------

# LCM: D14 -> D15
$ws.Range("C15").Value = "LCM"
$ws.Range("C15").HorizontalAlignment = -4152
$ws.Range("D15").Formula = "=LCM(24, 36)"

# GCD: D13 -> D14
$ws.Range("C14").Value = "GCD"
$ws.Range("C14").HorizontalAlignment = -4152
$ws.Range("D14").Formula = "=GCD(56, 21)"

# old GCD/LCM source cells are now empty
$ws.Range("C13:D13").Clear()

# PERMUT: D11 -> D12
$ws.Range("C12").Value = "PERMUT"
$ws.Range("C12").HorizontalAlignment = -4152
$ws.Range("D12").Formula = "=PERMUT(6,3)"

# COMBIN: D10 -> D11
$ws.Range("C11").Value = "COMBIN"
$ws.Range("C11").HorizontalAlignment = -4152
$ws.Range("D11").Formula = "=COMBIN(6,3)"

# old COMBIN source cells are now empty
$ws.Range("C10:D10").Clear()

# New CEIL row
$ws.Range("C9").Value = "CEIL"
$ws.Range("C9").HorizontalAlignment = -4152
$ws.Range("D9").Formula = "=CEILING(18.7, 1)"

# New PI row (note: shared-string insertion order matters, so PI/MOD/FACT
# are written in the same order the author typed them: PI, MOD, then FACT)
$ws.Range("C18").Value = "PI"
$ws.Range("C18").HorizontalAlignment = -4152
$ws.Range("D18").Formula = "=PI()"

# New MOD row
$ws.Range("C19").Value = "MOD"
$ws.Range("C19").HorizontalAlignment = -4152
$ws.Range("D19").Formula = "=MOD(5,3)"

# New FACT row
$ws.Range("C17").Value = "FACT"
$ws.Range("C17").HorizontalAlignment = -4152
$ws.Range("D17").Formula = "=FACT(6)"

# Widen column D to fit the new content
$ws.Columns("D").ColumnWidth = 13

# ---------------------------------------------------------------------------
# 2. Update the named ranges that describe these cells.
# ---------------------------------------------------------------------------
$wb.Names.Item("CombinVal").RefersTo = "=Sheet1!`$D`$11"
$wb.Names.Item("PermutVal").RefersTo = "=Sheet1!`$D`$12"
$wb.Names.Item("GcdVal").RefersTo = "=Sheet1!`$D`$14"
$wb.Names.Item("LcmVal").RefersTo = "=Sheet1!`$D`$15"

$wb.Names.Add("CeilVal", "=Sheet1!`$D`$9")
$wb.Names.Add("FACT", "=Sheet1!`$C`$17")
$wb.Names.Add("FactVal", "=Sheet1!`$D`$17")
$wb.Names.Add("PiVal", "=Sheet1!`$D`$18")
$wb.Names.Add("ModVal", "=Sheet1!`$D`$19")

# ---------------------------------------------------------------------------
# 3. Cosmetic: match the author's final selection.
# ---------------------------------------------------------------------------
$ws.Range("C18").Select()

Write-Output "edit applied"
